$d = $word.ActiveDocument

# 1. Header cell: "PHÒNG HÀNH CHÍNH" -> "PHÒNG ĐÀO TẠO" (paragraph 2, unique occurrence)
$r = $d.Paragraphs(2).Range
$r.Find.Execute("PHÒNG HÀNH CHÍNH", $true, $false, $false, $false, $false, `
                $true, 0, $false, "PHÒNG ĐÀO TẠO", 2)

# 2. Date: "05/11/2024" -> "14/10/2024" (unique occurrence in the whole document)
$r = $d.Content
$r.Find.Execute("05/11/2024", $true, $false, $false, $false, $false, `
                $true, 0, $false, "14/10/2024", 2)

# 3. "Ngày ..., Phòng Hành chính nhận được Công văn sau:" -> replace "Phòng Hành chính" with "Phòng Đào tạo"
$r = $d.Paragraphs(12).Range
$r.Find.Execute("Phòng Hành chính", $true, $false, $false, $false, $false, `
                $true, 0, $false, "Phòng Đào tạo", 2)

# 4. "Ý kiến của Phòng Hành chính: ${senderComment}" -> replace "Phòng Hành chính" with "Phòng Đào tạo"
$r = $d.Paragraphs(15).Range
$r.Find.Execute("Phòng Hành chính", $true, $false, $false, $false, $false, `
                $true, 0, $false, "Phòng Đào tạo", 2)
